$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "ELT-3A-Motores de aplicação"
$ws.Range("C3").Value = "ELT-3A-Motores de aplicação"
$ws.Range("C6").Value = "MCT-3A-Motores de aplicação"
$ws.Range("C7").Value = "-"
